$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking strings
# (e.g. "1.002", "316.62") are preserved verbatim as text, matching
# the original inlineStr cell type instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.524.33"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("E3").Value = "  +0.07%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "316.62"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.5195"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "0.3887"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "0.08388"
$ws.Range("E9").Value = "  +8.90%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "41.91"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "6.426"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("D13").Value = "21.13"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "1.002"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "7.520"
$ws.Range("E15").Value = "  -0.42%  "
$ws.Range("D16").Value = "1.821.87"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "0.00001128"
$ws.Range("E17").Value = "  +4.41%  "
$ws.Range("D18").Value = "93.23"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").Value = "0.06593"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").Value = "17.79"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  -0.11%  "
$ws.Range("D22").Value = "6.078"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "28.555.25"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +2.56%  "
$ws.Range("D25").Value = "2.278"
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "21.11"
$ws.Range("E26").Value = "  +2.29%  "
$ws.Range("D27").Value = "159.42"
$ws.Range("E27").Value = "  +1.78%  "
$ws.Range("D28").Value = "2.031.81"
$ws.Range("E28").Value = "  -0.28%  "
$ws.Range("D29").Value = "2.406"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "125.79"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").Value = "0.1095"
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("D32").Value = "1.100"
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("D33").Value = "5.738"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "0.07459"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").Value = "3.663"
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").Value = "0.2219"
$ws.Range("E36").Value = "  -1.25%  "
$ws.Range("D37").Value = "0.02366"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "5.225"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").Value = "8.811"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").Value = "11.48"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("D41").Value = "0.6331"
$ws.Range("E41").Value = "  +1.26%  "
$ws.Range("D42").Value = "1.192"
$ws.Range("E42").Value = "  +0.58%  "
$ws.Range("D43").Value = "1.399"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "13.53"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "3.784"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").Value = "0.5974"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("D47").Value = "126.98"
$ws.Range("E47").Value = "  +1.79%  "
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "0.06984"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "74.57"

Write-Output "Applied crypto list update"
